# Update FFXIV Leve profit-tracking sheets (Malboro_Profits) with refreshed
# market-board price snapshots. Mirrors a scheduled-runner data sync: for each
# changed row we rewrite the live price columns (H/I/J/K/L) and the derived
# profit columns (M/N); a couple of rows also gain/lose a cell entirely where
# the upstream snapshot no longer has an HQ (or NQ) price to report.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 511.66666
$ws.Range("I9").Value = 83
$ws.Range("K9").Value = 83
$ws.Range("M9").Value = 86
$ws.Range("H43").Value = 3899.077
$ws.Range("J43").Value = 5865.1665
$ws.Range("L43").Value = 5865.1665
$ws.Range("N43").Value = -6003.1665
$ws.Range("H112").Value = 6121.9546
$ws.Range("J112").Value = 3631.353
$ws.Range("L112").Value = 10894.059
$ws.Range("N112").Value = -13110.059
$ws.Range("H131").Value = 386401
$ws.Range("I131").Value = 386401
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 1159203
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = -1154163
$ws.Range("N131").ClearContents()
$ws.Range("H136").Value = 95053.57000000001
$ws.Range("J136").Value = 95053.57000000001
$ws.Range("L136").Value = 95053.57000000001
$ws.Range("N136").Value = -105253.57

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3496.3684
$ws.Range("I45").Value = 2857.3
$ws.Range("J45").Value = 4206.4443
$ws.Range("K45").Value = 2857.3
$ws.Range("L45").Value = 4206.4443
$ws.Range("M45").Value = -2480.3
$ws.Range("N45").Value = -4960.4443
$ws.Range("H74").Value = 31434
$ws.Range("I74").Value = 1925.0769
$ws.Range("K74").Value = 1925.0769
$ws.Range("M74").Value = -1051.0769
$ws.Range("H77").Value = 31434
$ws.Range("I77").Value = 1925.0769
$ws.Range("K77").Value = 9625.3845
$ws.Range("M77").Value = -5257.3845
$ws.Range("H97").Value = 1240.4
$ws.Range("J97").Value = 2166.6667
$ws.Range("L97").Value = 2166.6667
$ws.Range("N97").Value = -3158.6667
$ws.Range("H114").Value = 25000
$ws.Range("J114").Value = 25000
$ws.Range("L114").Value = 25000
$ws.Range("N114").Value = -33678
$ws.Range("H122").Value = 3514.2632
$ws.Range("I122").Value = 1779.1818
$ws.Range("J122").Value = 5900
$ws.Range("K122").Value = 5337.5454
$ws.Range("L122").Value = 17700
$ws.Range("M122").Value = -2887.5454
$ws.Range("N122").Value = -22600
$ws.Range("H132").Value = 3400122
$ws.Range("I132").Value = 1347.1515
$ws.Range("J132").Value = 22093384
$ws.Range("K132").Value = 4041.4545
$ws.Range("L132").Value = 66280152
$ws.Range("M132").Value = -1511.4545
$ws.Range("N132").Value = -66285212

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 37874.227
$ws.Range("I20").Value = 19948
$ws.Range("K20").Value = 19948
$ws.Range("M20").Value = -19701
$ws.Range("H94").Value = 1608.4062
$ws.Range("I94").Value = 1325.0435
$ws.Range("K94").Value = 1325.0435
$ws.Range("M94").Value = -874.0435
$ws.Range("H99").Value = 30917.066
$ws.Range("I99").Value = 31523.363
$ws.Range("K99").Value = 31523.363
$ws.Range("M99").Value = -30025.363

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5786.8125
$ws.Range("I16").Value = 2276.7778
$ws.Range("J16").Value = 10299.714
$ws.Range("K16").Value = 2276.7778
$ws.Range("L16").Value = 10299.714
$ws.Range("M16").Value = -1989.7778
$ws.Range("N16").Value = -10873.714
$ws.Range("H22").Value = 1451.1
$ws.Range("I22").Value = 450
$ws.Range("K22").Value = 450
$ws.Range("M22").Value = -100
$ws.Range("H31").Value = 15194.782
$ws.Range("I31").Value = 1634.909
$ws.Range("K31").Value = 1634.909
$ws.Range("M31").Value = -1339.909
$ws.Range("H34").Value = 15194.782
$ws.Range("I34").Value = 1634.909
$ws.Range("K34").Value = 1634.909
$ws.Range("M34").Value = -1432.909
$ws.Range("H105").Value = 7275.278
$ws.Range("I105").Value = 13808.5
$ws.Range("K105").Value = 13808.5
$ws.Range("M105").Value = -12061.5
$ws.Range("H113").Value = 5786.8125
$ws.Range("I113").Value = 2276.7778
$ws.Range("J113").Value = 10299.714
$ws.Range("K113").Value = 2276.7778
$ws.Range("L113").Value = 10299.714
$ws.Range("M113").Value = -106.7777999999998
$ws.Range("N113").Value = -14639.714
$ws.Range("H140").Value = 126000
$ws.Range("J140").Value = 126000
$ws.Range("L140").Value = 126000
$ws.Range("N140").Value = -136360

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 193266.33
$ws.Range("I9").Value = 257524.28
$ws.Range("K9").Value = 772572.84
$ws.Range("M9").Value = -772348.84

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()
$ws.Range("H70").Value = 4828.263
$ws.Range("I70").Value = 4724.143
$ws.Range("J70").Value = 4889
$ws.Range("K70").Value = 4724.143
$ws.Range("L70").Value = 4889
$ws.Range("M70").Value = -4454.143
$ws.Range("N70").Value = -5429
$ws.Range("H73").Value = 4828.263
$ws.Range("I73").Value = 4724.143
$ws.Range("J73").Value = 4889
$ws.Range("K73").Value = 4724.143
$ws.Range("L73").Value = 4889
$ws.Range("M73").Value = -3788.143
$ws.Range("N73").Value = -6761
$ws.Range("H113").Value = 3561.6
$ws.Range("I113").Value = 3303.6667
$ws.Range("K113").Value = 3303.6667
$ws.Range("M113").Value = -1133.6667
$ws.Range("H120").Value = 30000
$ws.Range("J120").Value = 30000
$ws.Range("L120").Value = 30000
$ws.Range("N120").Value = -39676
$ws.Range("H141").Value = 90610.5
$ws.Range("J141").Value = 90610.5
$ws.Range("L141").Value = 90610.5
$ws.Range("N141").Value = -100970.5

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2065
$ws.Range("I22").Value = 1233.3334
$ws.Range("K22").Value = 1233.3334
$ws.Range("M22").Value = -938.3334
$ws.Range("H27").Value = 2065
$ws.Range("I27").Value = 1233.3334
$ws.Range("K27").Value = 1233.3334
$ws.Range("M27").Value = -1126.3334
$ws.Range("H122").Value = 5339.2646
$ws.Range("I122").Value = 4182.4287
$ws.Range("J122").Value = 6149.05
$ws.Range("K122").Value = 12547.2861
$ws.Range("L122").Value = 18447.15
$ws.Range("M122").Value = -10097.2861
$ws.Range("N122").Value = -23347.15

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 113.2
$ws.Range("I81").Value = 113.2
$ws.Range("K81").Value = 226.4
$ws.Range("M81").Value = 834.6
$ws.Range("H84").Value = 113.2
$ws.Range("I84").Value = 113.2
$ws.Range("K84").Value = 1132
$ws.Range("M84").Value = 4172
$ws.Range("H122").Value = 3674.8286
$ws.Range("I122").Value = 2273.8076
$ws.Range("K122").Value = 6821.4228
$ws.Range("M122").Value = -4371.4228
$ws.Range("H138").Value = 49250
$ws.Range("J138").Value = 49250
$ws.Range("L138").Value = 49250
$ws.Range("N138").Value = -59530

